# DPLKAKT008-021 -> update test-data cell contents (username 32382 -> 31160)
# and swap displayed text between F2 / O2, plus move selection to F3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = "Username : 31160;`nPassword : bni1234;`nTgl. Market : 21/01/2023;`nFile Excel : 21012023HargaPasarFixedIncome.xlsx"
$ws.Range("G2").Value = 31160
$ws.Range("O2").Value = "21012023HargaPasarFixedIncome.xlsx"

$ws.Range("F3").Select()
